$wb = $excel.ActiveWorkbook

# 1. Rename existing sheet "Sheet1" -> "Search"
$search = $wb.Worksheets.Item(1)
$search.Name = "Search"

# 2. Add a new worksheet named "Hotel" right after "Search"
$hotel = $wb.Worksheets.Add($null, $search)
$hotel.Name = "Hotel"

# 3. Adjust the window height of the workbook view
$excel.ActiveWindow.Height = 9000

# 4. Update style of the label cells (row 1 + A2) on the Search sheet:
#    center horizontal + vertical alignment
$headerRange = $search.Range("A1:C1")
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

$a2Range = $search.Range("A2")
$a2Range.HorizontalAlignment = -4108
$a2Range.VerticalAlignment = -4108

# 5. Populate Hotel sheet content (centered horizontally, like the other
#    "center" style already used for the date cells)
$hotel.Range("A1").Value = "Hotel"
$hotel.Range("A2").Value = "Tolip"
$hotel.Range("A1:A2").HorizontalAlignment = -4108
